# Updates cryptos list cell values per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.406.75"
$ws.Range("E2").Value = "  +3.98%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.804.40"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.05"
$ws.Range("E5").Value = "  +0.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5525"
$ws.Range("E7").Value = "  +5.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3856"
$ws.Range("E8").Value = "  +6.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07589"
$ws.Range("E9").Value = "  +2.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.49"
$ws.Range("E10").Value = "  -0.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.127"
$ws.Range("E11").Value = "  +3.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.16%  "

$ws.Range("E13").Value = "  +3.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.192"
$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.406"
$ws.Range("E15").Value = "  +6.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.809.11"
$ws.Range("E16").Value = "  +1.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.12"
$ws.Range("E17").Value = "  +4.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001071"
$ws.Range("E18").Value = "  +2.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06444"
$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("E20").Value = "  -0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.33"
$ws.Range("E21").Value = "  +3.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.976"
$ws.Range("E22").Value = "  +2.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.425.95"
$ws.Range("E23").Value = "  +3.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.132"
$ws.Range("E25").Value = "  +2.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.27"
$ws.Range("E26").Value = "  +3.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.68"
$ws.Range("E27").Value = "  +2.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.409"
$ws.Range("E28").Value = "  +2.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.013.91"
$ws.Range("E29").Value = "  +1.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.97"
$ws.Range("E30").Value = "  +2.15%  "

$ws.Range("E31").Value = "  +5.21%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1020"
$ws.Range("E32").Value = "  +4.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.749"
$ws.Range("E33").Value = "  +3.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.683"
$ws.Range("E34").Value = "  +1.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.2314"
$ws.Range("E35").Value = "  +14.21%  "

$ws.Range("E36").Value = "  +8.12%  "

$ws.Range("E37").Value = "  +4.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.817"
$ws.Range("E38").Value = "  +8.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.136"
$ws.Range("E39").Value = "  +5.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.63"
$ws.Range("E40").Value = "  +3.75%  "

$ws.Range("E41").Value = "  +4.41%  "

$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9990"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.160"
$ws.Range("E43").Value = "  +1.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.386"
$ws.Range("E44").Value = "  -2.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.56"
$ws.Range("E45").Value = "  +2.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5989"
$ws.Range("E46").Value = "  +3.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.681"
$ws.Range("E47").Value = "  +1.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.27"
$ws.Range("E48").Value = "  +4.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.986"
$ws.Range("E49").Value = "  +5.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.148"
$ws.Range("E50").Value = "  +3.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06901"
$ws.Range("E51").Value = "  +2.83%  "
